$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibition listing) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 8696
$ws1.Range("F3").Value = 90
$ws1.Range("F4").Value = 233
$ws1.Range("F5").Value = 95
$ws1.Range("F6").Value = 1399
$ws1.Range("F7").Value = 1377
$ws1.Range("F8").Value = 232
$ws1.Range("F9").Value = 37
$ws1.Range("F10").Value = 271
$ws1.Range("F11").Value = 79

# --- Sheet "全部类型" (all types listing) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 8696
$ws4.Range("F3").Value = 90
$ws4.Range("F4").Value = 233
$ws4.Range("F5").Value = 95
$ws4.Range("F6").Value = 1399
$ws4.Range("F7").Value = 1377
$ws4.Range("F8").Value = 232
$ws4.Range("F10").Value = 37
$ws4.Range("F11").Value = 271
$ws4.Range("F12").Value = 79
